# Scheduled-runner refresh of price/profit figures on the Sheets tabs
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Updates currentAveragePrice* /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H:N)
# for the rows whose market data changed, including a couple of rows
# that gained or lost a LeveProfit cell entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 234.6
$ws.Range("I33").Value = 143.5
$ws.Range("K33").Value = 143.5
$ws.Range("M33").Value = 85.5
# Row 92
$ws.Range("H92").Value = 1819.1666
$ws.Range("I92").Value = 1819.1666
$ws.Range("K92").Value = 1819.1666
$ws.Range("M92").Value = -571.1666
# Row 98
$ws.Range("H98").Value = 898.5714
$ws.Range("I98").Value = 898.5714
$ws.Range("K98").Value = 898.5714
$ws.Range("M98").Value = 599.4286
# Row 107
$ws.Range("H107").Value = 285.85715
$ws.Range("I107").Value = 165.83333
$ws.Range("K107").Value = 165.83333
$ws.Range("M107").Value = 1754.16667
# Row 113
$ws.Range("H113").Value = 4068.125
$ws.Range("I113").Value = 2819.5
$ws.Range("J113").Value = 6149.1665
$ws.Range("K113").Value = 2819.5
$ws.Range("L113").Value = 6149.1665
$ws.Range("M113").Value = 434.5
$ws.Range("N113").Value = -12657.1665
# Row 116
$ws.Range("H116").Value = 6499
$ws.Range("I116").Value = 5665.1665
$ws.Range("J116").Value = 9000.5
$ws.Range("K116").Value = 5665.1665
$ws.Range("L116").Value = 9000.5
$ws.Range("M116").Value = -2223.1665
$ws.Range("N116").Value = -15884.5
# Row 122
$ws.Range("H122").Value = 898.5714
$ws.Range("I122").Value = 898.5714
$ws.Range("K122").Value = 2695.7142
$ws.Range("M122").Value = -245.7142000000003

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4799.222
$ws.Range("I61").Value = 3657.1667
$ws.Range("J61").Value = 7083.3335
$ws.Range("K61").Value = 3657.1667
$ws.Range("L61").Value = 7083.3335
$ws.Range("M61").Value = -3445.1667
$ws.Range("N61").Value = -7507.3335
# Row 74
$ws.Range("H74").Value = 1807.8235
$ws.Range("I74").Value = 1845.8125
$ws.Range("K74").Value = 1845.8125
$ws.Range("M74").Value = -971.8125
# Row 77
$ws.Range("H77").Value = 1807.8235
$ws.Range("I77").Value = 1845.8125
$ws.Range("K77").Value = 9229.0625
$ws.Range("M77").Value = -4861.0625
# Row 122
$ws.Range("H122").Value = 3291
$ws.Range("I122").Value = 3291
$ws.Range("K122").Value = 9873
$ws.Range("M122").Value = -7423
# Row 132
$ws.Range("H132").Value = 2688.6667
$ws.Range("I132").Value = 2751.4375
$ws.Range("J132").Value = 2186.5
$ws.Range("K132").Value = 8254.3125
$ws.Range("L132").Value = 6559.5
$ws.Range("M132").Value = -5724.3125
$ws.Range("N132").Value = -11619.5
# Row 136
$ws.Range("H136").Value = 4799.222
$ws.Range("I136").Value = 3657.1667
$ws.Range("J136").Value = 7083.3335
$ws.Range("K136").Value = 10971.5001
$ws.Range("L136").Value = 21250.0005
$ws.Range("M136").Value = -8421.500100000001
$ws.Range("N136").Value = -26350.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1647.7778
$ws.Range("I105").Value = 1628
$ws.Range("J105").Value = 1672.5
$ws.Range("K105").Value = 1628
$ws.Range("L105").Value = 1672.5
$ws.Range("M105").Value = 119
$ws.Range("N105").Value = -5166.5
# Row 134
$ws.Range("H134").Value = 4541.2
$ws.Range("I134").Value = 4601.3335
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 13804.0005
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -11269.0005
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4366.625
$ws.Range("I31").Value = 2057.7334
$ws.Range("K31").Value = 2057.7334
$ws.Range("M31").Value = -1762.7334
# Row 34
$ws.Range("H34").Value = 4366.625
$ws.Range("I34").Value = 2057.7334
$ws.Range("K34").Value = 2057.7334
$ws.Range("M34").Value = -1855.7334
# Row 132
$ws.Range("H132").Value = 1456.3334
$ws.Range("I132").Value = 1148.4
$ws.Range("K132").Value = 3445.2
$ws.Range("M132").Value = -915.2000000000003

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3700.2
$ws.Range("J68").Value = 4125.75
$ws.Range("L68").Value = 12377.25
$ws.Range("N68").Value = -13999.25
# Row 71
$ws.Range("H71").Value = 3700.2
$ws.Range("J71").Value = 4125.75
$ws.Range("L71").Value = 37131.75
$ws.Range("N71").Value = -45243.75
# Row 136
$ws.Range("H136").Value = 7675
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5166.6665
$ws.Range("I70").Value = 4250
$ws.Range("K70").Value = 4250
$ws.Range("M70").Value = -3980
# Row 73
$ws.Range("H73").Value = 5166.6665
$ws.Range("I73").Value = 4250
$ws.Range("K73").Value = 4250
$ws.Range("M73").Value = -3314
# Row 132
$ws.Range("H132").Value = 3723.3635
$ws.Range("I132").Value = 3795.7
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11387.1
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -8857.099999999999
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 832
$ws.Range("I16").Value = 832
$ws.Range("K16").Value = 832
$ws.Range("M16").Value = -662
# Row 18
$ws.Range("H18").Value = 8722.166999999999
$ws.Range("I18").Value = 9466.666999999999
$ws.Range("J18").Value = 4999.6665
$ws.Range("K18").Value = 9466.666999999999
$ws.Range("L18").Value = 4999.6665
$ws.Range("M18").Value = -9294.666999999999
$ws.Range("N18").Value = -5343.6665
# Row 20
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
# Row 22
$ws.Range("H22").Value = 4424.3335
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4424.3335
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4424.3335
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -5014.3335
# Row 27
$ws.Range("H27").Value = 4424.3335
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4424.3335
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4424.3335
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4638.3335
# Row 61
$ws.Range("H61").Value = 3842.7896
$ws.Range("I61").Value = 2951.2666
$ws.Range("J61").Value = 7186
$ws.Range("K61").Value = 2951.2666
$ws.Range("L61").Value = 7186
$ws.Range("M61").Value = -2749.2666
$ws.Range("N61").Value = -7590
# Row 63
$ws.Range("H63").Value = 44442.5
$ws.Range("I63").Value = 44442.5
$ws.Range("K63").Value = 44442.5
$ws.Range("M63").Value = -43693.5
# Row 66
$ws.Range("H66").Value = 44442.5
$ws.Range("I66").Value = 44442.5
$ws.Range("K66").Value = 133327.5
$ws.Range("M66").Value = -129583.5
# Row 98
$ws.Range("H98").Value = 50947
$ws.Range("J98").Value = 50947
$ws.Range("L98").Value = 50947
$ws.Range("N98").Value = -56937
# Row 113
$ws.Range("H113").Value = 3842.7896
$ws.Range("I113").Value = 2951.2666
$ws.Range("J113").Value = 7186
$ws.Range("K113").Value = 2951.2666
$ws.Range("L113").Value = 7186
$ws.Range("M113").Value = -781.2665999999999
$ws.Range("N113").Value = -11526
# Row 116
$ws.Range("H116").Value = 349999.5
$ws.Range("J116").Value = 349999.5
$ws.Range("L116").Value = 349999.5
$ws.Range("N116").Value = -359177.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 4656.6665
$ws.Range("I126").Value = 2521.4285
$ws.Range("K126").Value = 7564.2855
$ws.Range("M126").Value = -5094.2855
# Row 132
$ws.Range("H132").Value = 1798.4546
$ws.Range("I132").Value = 1808.6666
$ws.Range("K132").Value = 5425.9998
$ws.Range("M132").Value = -2895.9998
